$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.115.07"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.70"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.93"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.496"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.796.60"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.638.35"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.120.18"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.74"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.07"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.278.37"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.609"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -7.41%  "
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.56"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.31"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.709.25"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.08"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.74"
$ws.Range("E51").Value = "  +10.87%  "
